$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New list of student names (replacing the old placeholder names),
# mapped row-by-row onto B2:B17 in the same order as the shared-string table.
$names = @(
    "ARTHUR MAZZARDO NAUE",
    "CARLOS ALBERTO DE SOUZA",
    "DANIEL DIFENTTHAELER SANTOS",
    "GABRIEL LAZZARI OLIVEIRA",
    "GABRIEL MENDES ALVES ORTIZ PAULO",
    "HIGOR RODRIGUES ESTEVÃO",
    "JOAO PEDRO CARDOSO PERFEITO",
    "LUIZ FELIPE SCHALATA PACHECO",
    "MARIA EDUARDA DE OLIVEIRA ALVES",
    "MARILIA STEFENON RODRIGUES",
    "MICHEL DAVID DE SOUZA",
    "OSWALDO GABRIEL CARDOSO CORRÊA",
    "PATRICIA VOIGT",
    "PEDRO JAREMCZUK ZANONI SILVEIRA",
    "RODRIGO RAMOS DE OLIVEIRA",
    "VICTOR LEONARDO FAGUNDES DOS SANTOS"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $names[$i]
    $cell.Style = "Normal"
}

# Update the saved selection on the sheet to I16.
$ws.Range("I16").Select()
